$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two entirely-dropped rows (top-down order preserved by deleting
#     the later row first so the earlier row number stays valid) ---
# Row 26 = "RM 232", Row 28 = "SC 92" (original row numbers, before any shifting)
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# --- Per-cell value updates (row numbers below are AFTER the two row deletions
#     above, i.e. final/target row numbers) ---

# D5 (RM 14): -14.4 -> missing
$ws.Range("D5").ClearContents()

# E6 (RM 21): missing -> -5.7
$ws.Range("E6").Value = -5.7

# D8 (RM 38): missing -> -13.9
$ws.Range("D8").Value = -13.9

# E11 (RM 58): missing -> -7.9
$ws.Range("E11").Value = -7.9

# D12 (RM 81): -14.1 -> missing
$ws.Range("D12").ClearContents()
# E12 (RM 81): -5.3 -> missing
$ws.Range("E12").ClearContents()

# D14 (RM 90): missing -> -13.1
$ws.Range("D14").Value = -13.1

# E17 (RM 116): -7.3 -> missing
$ws.Range("E17").ClearContents()

# D18 (RM 120): -15.2 -> missing
$ws.Range("D18").ClearContents()

# E25 (RM 145): missing -> -7.1
$ws.Range("E25").Value = -7.1

# B26 (SC 5): missing -> -20.2
$ws.Range("B26").Value = -20.2

# B27 (SC 101): -20.4 -> missing
$ws.Range("B27").ClearContents()

# E31 (SC 132): -8.1 -> missing
$ws.Range("E31").ClearContents()

# E32 (SC 193): -6.4 -> missing
$ws.Range("E32").ClearContents()

# C33 (SC 232): missing -> 10.4
$ws.Range("C33").Value = 10.4
